# Applies the Leve-profit recalculation update from the scheduled runner.
# For each touched row, write the recomputed currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ) and LeveProfit(NQ/HQ) columns (H:N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 5073.75
$ws.Range("I19").Value = 1273.875
$ws.Range("J19").Value = 8873.625
$ws.Range("K19").Value = 1273.875
$ws.Range("L19").Value = 8873.625
$ws.Range("M19").Value = -1098.875
$ws.Range("N19").Value = -9223.625

# Row 33
$ws.Range("H33").Value = 800.3913
$ws.Range("I33").Value = 389.86667
$ws.Range("J33").Value = 1570.125
$ws.Range("K33").Value = 389.86667
$ws.Range("L33").Value = 1570.125
$ws.Range("M33").Value = -160.86667
$ws.Range("N33").Value = -2028.125

# Row 62
$ws.Range("H62").Value = 7043
$ws.Range("I62").Value = 6069.4
$ws.Range("K62").Value = 6069.4
$ws.Range("M62").Value = -5445.4

# Row 65
$ws.Range("H65").Value = 7043
$ws.Range("I65").Value = 6069.4
$ws.Range("K65").Value = 30347
$ws.Range("M65").Value = -27227

# Row 88
$ws.Range("H88").Value = 2121.8572
$ws.Range("I88").Value = 3834.3333
$ws.Range("K88").Value = 3834.3333
$ws.Range("M88").Value = -3428.3333

# Row 91
$ws.Range("H91").Value = 2121.8572
$ws.Range("I91").Value = 3834.3333
$ws.Range("K91").Value = 3834.3333
$ws.Range("M91").Value = -2430.3333

# Row 113
$ws.Range("H113").Value = 6490.923
$ws.Range("I113").Value = 4932.4443
$ws.Range("K113").Value = 4932.4443
$ws.Range("M113").Value = -1678.4443

# Row 133
$ws.Range("H133").Value = 120000
$ws.Range("J133").Value = 120000
$ws.Range("L133").Value = 120000
$ws.Range("N133").Value = -130120

# Row 135
$ws.Range("H135").Value = 9260714
$ws.Range("I135").Value = 953.0454999999999
$ws.Range("J135").Value = 50003660
$ws.Range("K135").Value = 8577.4095
$ws.Range("L135").Value = 450032940
$ws.Range("M135").Value = -6042.4095
$ws.Range("N135").Value = -450038010

# Row 137
$ws.Range("H137").Value = 48784076
$ws.Range("I137").Value = 32261016
$ws.Range("J137").Value = 100005560
$ws.Range("K137").Value = 96783048
$ws.Range("L137").Value = 300016680
$ws.Range("M137").Value = -96780498
$ws.Range("N137").Value = -300021780

# Row 138
$ws.Range("H138").Value = 6812.316
$ws.Range("I138").Value = 2967.4
$ws.Range("J138").Value = 8185.5
$ws.Range("K138").Value = 8902.200000000001
$ws.Range("L138").Value = 24556.5
$ws.Range("M138").Value = -3762.200000000001
$ws.Range("N138").Value = -34836.5

# Row 141
$ws.Range("H141").Value = 1628.2916
$ws.Range("I141").Value = 1667.6957
$ws.Range("K141").Value = 5003.0871
$ws.Range("M141").Value = 176.9129000000003

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17249556
$ws.Range("I32").Value = 19236236
$ws.Range("J32").Value = 31664.5
$ws.Range("K32").Value = 19236236
$ws.Range("L32").Value = 31664.5
$ws.Range("M32").Value = -19235949
$ws.Range("N32").Value = -32238.5

# Row 35
$ws.Range("H35").Value = 1361.25
$ws.Range("I35").Value = 347.5
$ws.Range("J35").Value = 2375
$ws.Range("K35").Value = 347.5
$ws.Range("L35").Value = 2375
$ws.Range("M35").Value = 58.5
$ws.Range("N35").Value = -3187

# Row 45
$ws.Range("H45").Value = 2080.182
$ws.Range("I45").Value = 1347.8334
$ws.Range("J45").Value = 2959
$ws.Range("K45").Value = 1347.8334
$ws.Range("L45").Value = 2959
$ws.Range("M45").Value = -970.8334
$ws.Range("N45").Value = -3713

# Row 61
$ws.Range("H61").Value = 47623850
$ws.Range("I61").Value = 71431460
$ws.Range("K61").Value = 71431460
$ws.Range("M61").Value = -71431248

# Row 74
$ws.Range("H74").Value = 91012120
$ws.Range("I74").Value = 91012120
$ws.Range("K74").Value = 91012120
$ws.Range("M74").Value = -91011246

# Row 77
$ws.Range("H77").Value = 91012120
$ws.Range("I77").Value = 91012120
$ws.Range("K77").Value = 455060600
$ws.Range("M77").Value = -455056232

# Row 122
$ws.Range("H122").Value = 3481.1333
$ws.Range("I122").Value = 2324.7144
$ws.Range("J122").Value = 4493
$ws.Range("K122").Value = 6974.1432
$ws.Range("L122").Value = 13479
$ws.Range("M122").Value = -4524.1432
$ws.Range("N122").Value = -18379

# Row 136
$ws.Range("H136").Value = 47623850
$ws.Range("I136").Value = 71431460
$ws.Range("K136").Value = 214294380
$ws.Range("M136").Value = -214291830

$ws = $wb.Worksheets.Item("BSM")
# Row 132
$ws.Range("H132").Value = 110000
$ws.Range("J132").Value = 110000
$ws.Range("L132").Value = 110000
$ws.Range("N132").Value = -120120

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1971.8823
$ws.Range("I16").Value = 1998.9231
$ws.Range("K16").Value = 1998.9231
$ws.Range("M16").Value = -1711.9231

# Row 31
$ws.Range("H31").Value = 20837732
$ws.Range("I31").Value = 3191.879
$ws.Range("K31").Value = 3191.879
$ws.Range("M31").Value = -2896.879

# Row 34
$ws.Range("H34").Value = 20837732
$ws.Range("I34").Value = 3191.879
$ws.Range("K34").Value = 3191.879
$ws.Range("M34").Value = -2989.879

# Row 58
$ws.Range("H58").Value = 1513.5769
$ws.Range("I58").Value = 1474.12
$ws.Range("K58").Value = 1474.12
$ws.Range("M58").Value = -1271.12

# Row 113
$ws.Range("H113").Value = 1971.8823
$ws.Range("I113").Value = 1998.9231
$ws.Range("K113").Value = 1998.9231
$ws.Range("M113").Value = 171.0769

# Row 122
$ws.Range("H122").Value = 2580.2942
$ws.Range("I122").Value = 2430.182
$ws.Range("J122").Value = 2855.5
$ws.Range("K122").Value = 7290.545999999999
$ws.Range("L122").Value = 8566.5
$ws.Range("M122").Value = -4840.545999999999
$ws.Range("N122").Value = -13466.5

# Row 134
$ws.Range("H134").Value = 983.6977000000001
$ws.Range("I134").Value = 1029.5526
$ws.Range("J134").Value = 635.2
$ws.Range("K134").Value = 3088.6578
$ws.Range("L134").Value = 1905.6
$ws.Range("M134").Value = -553.6578
$ws.Range("N134").Value = -6975.6

# Row 136
$ws.Range("H136").Value = 1513.5769
$ws.Range("I136").Value = 1474.12
$ws.Range("K136").Value = 4422.36
$ws.Range("M136").Value = -1872.36

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1605.65
$ws.Range("J5").Value = 2948.9167
$ws.Range("L5").Value = 8846.750100000001
$ws.Range("N5").Value = -9070.750100000001

# Row 63
$ws.Range("H63").Value = 4936.6
$ws.Range("I63").Value = 5062.125
$ws.Range("J63").Value = 4434.5
$ws.Range("K63").Value = 15186.375
$ws.Range("L63").Value = 13303.5
$ws.Range("M63").Value = -14437.375
$ws.Range("N63").Value = -14801.5

# Row 66
$ws.Range("H66").Value = 4936.6
$ws.Range("I66").Value = 5062.125
$ws.Range("J66").Value = 4434.5
$ws.Range("K66").Value = 45559.125
$ws.Range("L66").Value = 39910.5
$ws.Range("M66").Value = -41815.125
$ws.Range("N66").Value = -47398.5

# Row 92
$ws.Range("H92").Value = 219.8
$ws.Range("I92").Value = 224.75
$ws.Range("J92").Value = 200
$ws.Range("K92").Value = 674.25
$ws.Range("L92").Value = 600
$ws.Range("M92").Value = 573.75
$ws.Range("N92").Value = -3096

# Row 129
$ws.Range("H129").Value = 3182.1
$ws.Range("I129").Value = 941
$ws.Range("K129").Value = 2823
$ws.Range("M129").Value = 2177

# Row 134
$ws.Range("H134").Value = 4459.5
$ws.Range("I134").Value = 1451.5333
$ws.Range("J134").Value = 19499.334
$ws.Range("K134").Value = 4354.5999
$ws.Range("L134").Value = 58498.00199999999
$ws.Range("M134").Value = 715.4000999999998
$ws.Range("N134").Value = -68638.00199999999

# Row 135
$ws.Range("H135").Value = 1605.65
$ws.Range("J135").Value = 2948.9167
$ws.Range("L135").Value = 26540.2503
$ws.Range("N135").Value = -31610.2503

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2777.4
$ws.Range("I122").Value = 2462.6667
$ws.Range("K122").Value = 7388.000100000001
$ws.Range("M122").Value = -4938.000100000001

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2342.3684
$ws.Range("I22").Value = 1231
$ws.Range("J22").Value = 3150.6365
$ws.Range("K22").Value = 1231
$ws.Range("L22").Value = 3150.6365
$ws.Range("M22").Value = -936
$ws.Range("N22").Value = -3740.6365

# Row 27
$ws.Range("H27").Value = 2342.3684
$ws.Range("I27").Value = 1231
$ws.Range("J27").Value = 3150.6365
$ws.Range("K27").Value = 1231
$ws.Range("L27").Value = 3150.6365
$ws.Range("M27").Value = -1124
$ws.Range("N27").Value = -3364.6365

# Row 68
$ws.Range("H68").Value = 5299.3335
$ws.Range("I68").Value = 3874.75
$ws.Range("J68").Value = 6439
$ws.Range("K68").Value = 3874.75
$ws.Range("L68").Value = 6439
$ws.Range("M68").Value = -3125.75
$ws.Range("N68").Value = -7937

# Row 71
$ws.Range("H71").Value = 5299.3335
$ws.Range("I71").Value = 3874.75
$ws.Range("J71").Value = 6439
$ws.Range("K71").Value = 19373.75
$ws.Range("L71").Value = 32195
$ws.Range("M71").Value = -15629.75
$ws.Range("N71").Value = -39683

# Row 122
$ws.Range("H122").Value = 5243.625
$ws.Range("I122").Value = 4662.25
$ws.Range("J122").Value = 5825
$ws.Range("K122").Value = 13986.75
$ws.Range("L122").Value = 17475
$ws.Range("M122").Value = -11536.75
$ws.Range("N122").Value = -22375

# Row 132
$ws.Range("H132").Value = 76927740
$ws.Range("I132").Value = 4806.8887
$ws.Range("J132").Value = 250004350
$ws.Range("K132").Value = 14420.6661
$ws.Range("L132").Value = 750013050
$ws.Range("M132").Value = -11890.6661
$ws.Range("N132").Value = -750018110

# Row 136
$ws.Range("H136").Value = 6202.16
$ws.Range("I136").Value = 5099.6665
$ws.Range("J136").Value = 11990.25
$ws.Range("K136").Value = 15298.9995
$ws.Range("L136").Value = 35970.75
$ws.Range("M136").Value = -12748.9995
$ws.Range("N136").Value = -41070.75

$ws = $wb.Worksheets.Item("WVR")
# Row 40
$ws.Range("H40").Value = 22747.5
$ws.Range("J40").Value = 22747.5
$ws.Range("L40").Value = 22747.5
$ws.Range("N40").Value = -23045.5

# Row 62
$ws.Range("H62").Value = 10299.6
$ws.Range("I62").Value = 10299.6
$ws.Range("K62").Value = 10299.6
$ws.Range("M62").Value = -9675.6

# Row 65
$ws.Range("H65").Value = 10299.6
$ws.Range("I65").Value = 10299.6
$ws.Range("K65").Value = 51498
$ws.Range("M65").Value = -48378

# Row 122
$ws.Range("H122").Value = 40044060
$ws.Range("I122").Value = 55615096
$ws.Range("K122").Value = 166845288
$ws.Range("M122").Value = -166842838
